$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.897803247373448
$ws.Range("B3").Value = 0.9214659685863874
$ws.Range("B4").Value = 0.8949152542372881
$ws.Range("B5").Value = 0.9079965606190885
